$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.057742643088169
$ws.Range("D2").Value = 1.055265475451534
$ws.Range("E2").Value = 1.063305949821053
$ws.Range("F2").Value = 1.072666294045031
$ws.Range("I2").Value = 1.048497942368824
$ws.Range("J2").Value = 1.062737480654799
$ws.Range("K2").Value = 1.058006408157197
$ws.Range("L2").Value = 1.066024946841664
$ws.Range("M2").Value = 1.075360202523923

# Row 3
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.058886898265979
$ws.Range("D3").Value = 1.056120190283747
$ws.Range("E3").Value = 1.064334219556289
$ws.Range("F3").Value = 1.073797905950379
$ws.Range("I3").Value = 1.048832845641868
$ws.Range("J3").Value = 1.063533389887483
$ws.Range("K3").Value = 1.058674424707046
$ws.Range("L3").Value = 1.066867666223849
$ws.Range("M3").Value = 1.076307826669175

# Row 4
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.059627357309382
$ws.Range("D4").Value = 1.056673222483377
$ws.Range("E4").Value = 1.064999939573588
$ws.Range("F4").Value = 1.074530671708417
$ws.Range("I4").Value = 1.049048396441891
$ws.Range("J4").Value = 1.064047873578672
$ws.Range("K4").Value = 1.059106006876108
$ws.Range("L4").Value = 1.067412712781813
$ws.Range("M4").Value = 1.076920948638579

# Row 5
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.059938658940513
$ws.Range("D5").Value = 1.05690571102451
$ws.Range("E5").Value = 1.065279894779086
$ws.Range("F5").Value = 1.074838855317227
$ws.Range("I5").Value = 1.049138737794471
$ws.Range("J5").Value = 1.064264037615964
$ws.Range("K5").Value = 1.059287283906703
$ws.Range("L5").Value = 1.067641790886076
$ws.Range("M5").Value = 1.077178692379643

# Row 6
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.059990928619758
$ws.Range("D6").Value = 1.056944746507281
$ws.Range("E6").Value = 1.065326905567316
$ws.Range("F6").Value = 1.074890608267277
$ws.Range("I6").Value = 1.049153890310421
$ws.Range("J6").Value = 1.064300325191845
$ws.Range("K6").Value = 1.059317711731517
$ws.Range("L6").Value = 1.067680250626315
$ws.Range("M6").Value = 1.077221967948304

# Row 7
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.059631516886277
$ws.Range("D7").Value = 1.056676329031115
$ws.Range("E7").Value = 1.065003680008629
$ws.Range("F7").Value = 1.074534789166095
$ws.Range("I7").Value = 1.049049604672667
$ws.Range("J7").Value = 1.064050762463594
$ws.Range("K7").Value = 1.059108429737676
$ws.Range("L7").Value = 1.067415773968074
$ws.Range("M7").Value = 1.076924392672789

# Row 8
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.058129339738296
$ws.Range("D8").Value = 1.055554335539
$ws.Range("E8").Value = 1.06365338378799
$ws.Range("F8").Value = 1.073048616527131
$ws.Range("I8").Value = 1.04861136348502
$ws.Range("J8").Value = 1.063006570308762
$ws.Range("K8").Value = 1.05823230609975
$ws.Range("L8").Value = 1.066309799527064
$ws.Range("M8").Value = 1.075680467914101

# Row 9
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.05548265348719
$ws.Range("D9").Value = 1.053577048780785
$ws.Range("E9").Value = 1.061276745052431
$ws.Range("F9").Value = 1.070433891671437
$ws.Range("I9").Value = 1.047830284881648
$ws.Range("J9").Value = 1.061162554552437
$ws.Range("K9").Value = 1.056683334283536
$ws.Range("L9").Value = 1.064359013159382
$ws.Range("M9").Value = 1.073488086938377

# Row 10
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.053718371838701
$ws.Range("D10").Value = 1.052258732058413
$ws.Range("E10").Value = 1.059694159173467
$ws.Range("F10").Value = 1.068693478064748
$ws.Range("I10").Value = 1.047303615240334
$ws.Range("J10").Value = 1.059930491954082
$ws.Range("K10").Value = 1.055647228444098
$ws.Range("L10").Value = 1.06305718383358
$ws.Range("M10").Value = 1.072026194208776

# Row 11
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.052954447290713
$ws.Range("D11").Value = 1.051687856299236
$ws.Range("E11").Value = 1.059009314651615
$ws.Range("F11").Value = 1.067940503135829
$ws.Range("I11").Value = 1.04707414830438
$ws.Range("J11").Value = 1.059396345447072
$ws.Range("K11").Value = 1.05519776031168
$ws.Range("L11").Value = 1.06249316358021
$ws.Range("M11").Value = 1.071393101094777

# Row 12
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.052670693653543
$ws.Range("D12").Value = 1.051475802099836
$ws.Range("E12").Value = 1.058754996417077
$ws.Range("F12").Value = 1.067660909628143
$ws.Range("I12").Value = 1.046988701185288
$ws.Range("J12").Value = 1.059197840681059
$ws.Range("K12").Value = 1.055030683194493
$ws.Range("L12").Value = 1.062283612872215
$ws.Range("M12").Value = 1.071157928985122

# Row 13
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.052731559702915
$ws.Range("D13").Value = 1.051521288697489
$ws.Range("E13").Value = 1.058809545674456
$ws.Range("F13").Value = 1.067720879099392
$ws.Range("I13").Value = 1.047007039519242
$ws.Range("J13").Value = 1.059240425114224
$ws.Range("K13").Value = 1.055066527448698
$ws.Range("L13").Value = 1.062328564411329
$ws.Range("M13").Value = 1.071208374804028

# Row 14
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.052930992093672
$ws.Range("D14").Value = 1.051670327946397
$ws.Range("E14").Value = 1.058988291333836
$ws.Range("F14").Value = 1.067917389916606
$ws.Range("I14").Value = 1.047067089562846
$ws.Range("J14").Value = 1.059379939003179
$ws.Range("K14").Value = 1.055183952209541
$ws.Range("L14").Value = 1.062475843044819
$ws.Range("M14").Value = 1.071373661965924

# Row 15
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.053053869194497
$ws.Range("D15").Value = 1.051762155210719
$ws.Range("E15").Value = 1.059098430837269
$ws.Range("F15").Value = 1.068038479275202
$ws.Range("I15").Value = 1.047104060161173
$ws.Range("J15").Value = 1.059465884984527
$ws.Range("K15").Value = 1.05525628495452
$ws.Range("L15").Value = 1.06256657983479
$ws.Range("M15").Value = 1.071475499084662

# Row 16
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.053769071335271
$ws.Range("D16").Value = 1.05229661843476
$ws.Range("E16").Value = 1.059739618979323
$ws.Range("F16").Value = 1.068743463872749
$ws.Range("I16").Value = 1.047318814363389
$ws.Range("J16").Value = 1.059965927648073
$ws.Range("K16").Value = 1.055677040702427
$ws.Range("L16").Value = 1.063094609227936
$ws.Range("M16").Value = 1.072068208700335

# Row 17
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.054217703524464
$ws.Range("D17").Value = 1.0526318633887
$ws.Range("E17").Value = 1.060141933005299
$ws.Range("F17").Value = 1.069185851848549
$ws.Range("I17").Value = 1.047453144875125
$ws.Range("J17").Value = 1.060279415395715
$ws.Range("K17").Value = 1.055940747902787
$ws.Range("L17").Value = 1.063425742233347
$ws.Range("M17").Value = 1.072439976941479

# Row 18
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.054479385201978
$ws.Range("D18").Value = 1.052827402641444
$ws.Range("E18").Value = 1.060376637063758
$ws.Range("F18").Value = 1.069443950640682
$ws.Range("I18").Value = 1.047531361028708
$ws.Range("J18").Value = 1.060462204305982
$ws.Range("K18").Value = 1.056094484034281
$ws.Range("L18").Value = 1.063618855710694
$ws.Range("M18").Value = 1.072656815005735

# Row 19
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.054568612286542
$ws.Range("D19").Value = 1.052894075890942
$ws.Range("E19").Value = 1.060456672065651
$ws.Range("F19").Value = 1.06953196608643
$ws.Range("I19").Value = 1.047558007559305
$ws.Range("J19").Value = 1.060524519898648
$ws.Range("K19").Value = 1.056146890553712
$ws.Range("L19").Value = 1.063684697205409
$ws.Range("M19").Value = 1.072730749830939

# Row 20
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.054169569326476
$ws.Range("D20").Value = 1.052595895136248
$ws.Range("E20").Value = 1.060098764270213
$ws.Range("F20").Value = 1.069138381461746
$ws.Range("I20").Value = 1.047438746609186
$ws.Range("J20").Value = 1.060245787666807
$ws.Range("K20").Value = 1.055912462873945
$ws.Range("L20").Value = 1.063390217984917
$ws.Range("M20").Value = 1.072400090565063

# Row 21
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.052872264181076
$ws.Range("D21").Value = 1.051626439755824
$ws.Range("E21").Value = 1.058935653442424
$ws.Range("F21").Value = 1.067859519763684
$ws.Range("I21").Value = 1.047049412200044
$ws.Range("J21").Value = 1.059338858376438
$ws.Range("K21").Value = 1.055149376987477
$ws.Range("L21").Value = 1.062432474499873
$ws.Range("M21").Value = 1.071324989383371

# Row 22
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.052056607653305
$ws.Range("D22").Value = 1.051016872655342
$ws.Range("E22").Value = 1.058204726759069
$ws.Range("F22").Value = 1.067055997467425
$ws.Range("I22").Value = 1.04680339042755
$ws.Range("J22").Value = 1.058768063064962
$ws.Range("K22").Value = 1.05466887331453
$ws.Range("L22").Value = 1.061830022480809
$ws.Range("M22").Value = 1.070648954321431

# Row 23
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.052489001999542
$ws.Range("D23").Value = 1.051340018787317
$ws.Range("E23").Value = 1.058592170155005
$ws.Range("F23").Value = 1.067481908016286
$ws.Range("I23").Value = 1.046933928054762
$ws.Range("J23").Value = 1.059070706935463
$ws.Range("K23").Value = 1.054923665880581
$ws.Range("L23").Value = 1.062149420453981
$ws.Range("M23").Value = 1.071007340706317

# Row 24
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.054191319091735
$ws.Range("D24").Value = 1.052612147651829
$ws.Range("E24").Value = 1.060118270236234
$ws.Range("F24").Value = 1.069159831096326
$ws.Range("I24").Value = 1.047445252987638
$ws.Range("J24").Value = 1.060260982785035
$ws.Range("K24").Value = 1.055925243907414
$ws.Range("L24").Value = 1.063406269958622
$ws.Range("M24").Value = 1.072418113526351

# Row 25
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.056166850383011
$ws.Range("D25").Value = 1.054088247130167
$ws.Range("E25").Value = 1.061890837888072
$ws.Range("F25").Value = 1.071109376884251
$ws.Range("I25").Value = 1.048033260820295
$ws.Range("J25").Value = 1.061639754651755
$ws.Range("K25").Value = 1.057084389469653
$ws.Range("L25").Value = 1.064863567092966
$ws.Range("M25").Value = 1.07405492244365
